$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the activation date "01/01/2012" -> "01/01/2023" ---
# The same shared text is used in B8/C8 ("Ativacao:") and (due to a pre-existing
# data artefact in the source) also in B13/C13 ("Programa resumido:").
# Force text (not a date serial) by setting the NumberFormat to Text first.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2023"

# --- Row 11 (under "Objectives:"): add the English objectives text ---
$ws.Range("B11").Font.Bold = $false
$ws.Range("B11").WrapText = $true
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").Value = "Provide student seminars on current topics in Physics, Technology and Engineering."

$ws.Range("C11").Font.Color = 255
$ws.Range("C11").WrapText = $true
$ws.Range("C11").VerticalAlignment = -4160
$ws.Range("C11").Value = "Provide student seminars on current topics in Physics, Technology and Engineering."

# --- Row 14 (under "Short syllabus:"): add the English short syllabus text ---
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").WrapText = $true
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("B14").Value = "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer."

$ws.Range("C14").Font.Color = 255
$ws.Range("C14").WrapText = $true
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C14").Value = "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer."

# --- Row 16 (under "Syllabus:"): add the English syllabus text ---
$ws.Range("B16").Font.Bold = $false
$ws.Range("B16").WrapText = $true
$ws.Range("B16").VerticalAlignment = -4160
$ws.Range("B16").Value = "Seminars followed by debates with professionals and undergraduate and graduate students on relevant and current topics in the areas of Physics, Technology and Engineering, ranging from basic research to the industrial and services segment."

$ws.Range("C16").Font.Color = 255
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4160
$ws.Range("C16").Value = "Seminars followed by debates with professionals and undergraduate and graduate students on relevant and current topics in the areas of Physics, Technology and Engineering, ranging from basic research to the industrial and services segment."
